$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("List4")
$ws.Range("A11:G32").FormulaArray = "=TRANSPOSE(A1:V7)"
